$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 6 (Subdomain), pushing the rest of the
# form fields (and the trailing blank row) down by one row.
$ws.Rows("6").Insert()

# Sponsor: drop the extra comma-separated users, keep just the admin.
$ws.Range("B3").Value = "mpfp-base-unilevel-business-admin"

# Subdomain value (now on row 7 after the insert).
$ws.Range("B7").Value = "rootcase20"

# Enrollment Package value (now on row 12 after the insert).
$ws.Range("B12").Value = "enrollment-package-9"
